# Scene 2B dialogue tweaks ("write some new for stephen")
#
# Four lines of dialogue are reworded. Each is looked up by its full,
# unique original sentence and swapped for the new wording with a single
# targeted Find/Replace so no other part of the document is touched.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Ah, well. At least the day’s over. I get up and stretch, ready to go home and relax…",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Well, at least the day’s over. I get up and stretch, ready to go home and relax…",
    2) | Out-Null

$d.Content.Find.Execute(
    "There it is. Well, might as well get it over with.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There it is. Might as well get it over with, I guess.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Asher (neutral curious): Sure, that’d be great. Your mom will be okay with it?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Asher (neutral curious): Sure, that’d be great. Will your mom be okay with it?",
    2) | Out-Null

$d.Content.Find.Execute(
    "Petra (neutral raised_eyebrow): You sure?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Petra (neutral curious): You sure?",
    2) | Out-Null
